$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 6
$ws.Range("B4").Value = 8

$ws.Range("B4").Select()
